# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the bc97d4fa... file row on both the
# zh-cn and de-de status sheets, reflecting a newly generated handoff report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for bc97d4fa-0c5c-4f17-9ed4-f9014c67d62b (row 4) -> column D
# "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-09 15:05:39"

# de-de sheet: same row/column for the same file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-09 15:05:52"
